$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values that changed with the new TPM numbers ---
$ws.Range("I2").Value = 0.9476581720434079
$ws.Range("J2").Value = 0.947658172043408
$ws.Range("M2").Value = 0.06202966666666667
$ws.Range("N2").Value = 0.186089
$ws.Range("Q2").Value = 0.1538496183404444
$ws.Range("R2").Value = 1.384646565064
$ws.Range("S2").Value = 0.9476581720434079
$ws.Range("T2").Value = 0.947658172043408

# --- Add the new row 3 (Resolving-Mac / Tac1 / Tacr2 / ECs) ---
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Tac1"
$ws.Range("C3").Value = "Tacr2"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1369916666666667
$ws.Range("H3").Value = 0.410975
$ws.Range("I3").Value = 0.0523418279565921
$ws.Range("J3").Value = 0.0523418279565921
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.06202966666666667
$ws.Range("N3").Value = 0.186089
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.008497547419444444
$ws.Range("R3").Value = 0.076477926775
$ws.Range("S3").Value = 0.0523418279565921
$ws.Range("T3").Value = 0.0523418279565921
